$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric cell updates
$ws.Range("A4").Value = 3791789.34
$ws.Range("B4").Value = 424.17
$ws.Range("C4").Value = 1486.98
$ws.Range("D4").Value = 338.49
$ws.Range("E4").Value = 4453089.86
$ws.Range("F4").Value = 563.36
$ws.Range("G4").Value = 1898.97
$ws.Range("H4").Value = 331.18
$ws.Range("A5").Value = 1137536.8
$ws.Range("B5").Value = 127.25
$ws.Range("C5").Value = 446.09
$ws.Range("D5").Value = 101.55
$ws.Range("E5").Value = 1335926.96
$ws.Range("F5").Value = 169.01
$ws.Range("G5").Value = 569.69
$ws.Range("H5").Value = 99.35
$ws.Range("C6").Value = 10421.05
$ws.Range("G6").Value = 10005.06
$ws.Range("H6").Value = 1744.88
$ws.Range("A7").Value = 2464663.07
$ws.Range("B7").Value = 275.71
$ws.Range("C7").Value = 966.53
$ws.Range("D7").Value = 220.02
$ws.Range("E7").Value = 2894508.41
$ws.Range("F7").Value = 366.18
$ws.Range("G7").Value = 1234.33
$ws.Range("H7").Value = 215.27
$ws.Range("A8").Value = 33967671.46
$ws.Range("B8").Value = 3799.76
$ws.Range("C8").Value = 13320.66
$ws.Range("D8").Value = 3032.24
$ws.Range("E8").Value = 32145394.34
$ws.Range("F8").Value = 4066.71
$ws.Range("G8").Value = 13708.06
$ws.Range("H8").Value = 2390.68
$ws.Range("A11").Value = 6318768.140000001
$ws.Range("B11").Value = 706.84
$ws.Range("C11").Value = 6194.87
$ws.Range("D11").Value = 355.03
$ws.Range("E11").Value = 6311627.84
$ws.Range("F11").Value = 798.48
$ws.Range("G11").Value = 7147.94
$ws.Range("H11").Value = 344.03
$ws.Range("A12").Value = 1895630.44
$ws.Range("B12").Value = 212.05
$ws.Range("C12").Value = 1858.46
$ws.Range("D12").Value = 106.51
$ws.Range("E12").Value = 1893488.35
$ws.Range("F12").Value = 239.54
$ws.Range("G12").Value = 2144.38
$ws.Range("H12").Value = 103.21
$ws.Range("C13").Value = 34577.48
$ws.Range("G13").Value = 35263.95
$ws.Range("H13").Value = 1697.25
$ws.Range("A14").Value = 4107199.29
$ws.Range("B14").Value = 459.45
$ws.Range("C14").Value = 4026.67
$ws.Range("D14").Value = 230.77
$ws.Range("E14").Value = 4102558.09
$ws.Range("F14").Value = 519.01
$ws.Range("G14").Value = 4646.16
$ws.Range("H14").Value = 223.62
$ws.Range("A15").Value = 47590624.49
$ws.Range("B15").Value = 5323.68
$ws.Range("C15").Value = 46657.47
$ws.Range("D15").Value = 2673.96
$ws.Range("E15").Value = 43445739.08
$ws.Range("F15").Value = 5496.31
$ws.Range("G15").Value = 49202.42
$ws.Range("H15").Value = 2368.1
$ws.Range("A17").Value = 15754862.96
$ws.Range("B17").Value = 1762.4
$ws.Range("C17").Value = 350.77
$ws.Range("E17").Value = 15857611.35
$ws.Range("F17").Value = 2006.14
$ws.Range("G17").Value = 339.09
$ws.Range("A18").Value = 4726458.89
$ws.Range("B18").Value = 528.72
$ws.Range("C18").Value = 105.23
$ws.Range("E18").Value = 4757283.41
$ws.Range("F18").Value = 601.84
$ws.Range("G18").Value = 101.73
$ws.Range("G19").Value = 377.37
$ws.Range("A20").Value = 30069062.24
$ws.Range("B20").Value = 3363.65
$ws.Range("C20").Value = 669.47
$ws.Range("E20").Value = 26495385.08
$ws.Range("F20").Value = 3351.92
$ws.Range("G20").Value = 566.56
$ws.Range("A21").Value = 10240660.93
$ws.Range("B21").Value = 1145.56
$ws.Range("C21").Value = 228
$ws.Range("E21").Value = 10307447.38
$ws.Range("F21").Value = 1303.99
$ws.Range("G21").Value = 220.41
$ws.Range("A22").Value = 80782201.38
$ws.Range("B22").Value = 9036.63
$ws.Range("C22").Value = 1798.57
$ws.Range("E22").Value = 75065366.19
$ws.Range("F22").Value = 9496.5
$ws.Range("G22").Value = 1605.15

# Text cell updates (labels)
$ws.Range("A9").Value = "ТО-1 после рек"
$ws.Range("E9").Value = "ТО-1 до рек"
$ws.Range("A16").Value = "ТО-2 после рек"
$ws.Range("E16").Value = "ТО-2 до рек"
$ws.Range("A23").Value = "ТР после рек"
$ws.Range("E23").Value = "ТР до рек"
